$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the redundant repeated header/value cells in H1:AU2 (keep styles, just remove content)
$ws.Range("H1:AU2").ClearContents()

# Add note about new file-writing step to D5
$ws.Range("D5").Value = "Writes new file without errors but keeps combined fluency.  Moves combinedflunecy.csv**   "

# Update the view: scroll so column R is the left-most visible column,
# and select H1:T2 (active cell lands on the anchor of the selection)
$excel.ActiveWindow.ScrollColumn = $ws.Range("R1").Column
$ws.Range("H1:T2").Select()
